$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (LeetCode #1193 - Monthly Transactions I)
$ws.Range("A16").Value = "1193. Monthly Transactions I"
$ws.Range("B16").Value = "Medium"
$ws.Range("E16").Value = "https://leetcode.com/problems/monthly-transactions-i/solutions/3871182/100-easy-fast-clean-solution/?envType=study-plan-v2&envId=top-sql-50 "
$ws.Range("D16").Value = "For grouping by month, use TO_CHAR(trans_date, 'YYYY-MM') AS month (or DATE_FORMAT). Use sum(case) for the counts. Group by TO_CHAR(trans_date, 'YYYY-MM'), country. A good thought process is to work column by column in the output table and treat each column as its own query, and work towards the solution."
$ws.Range("C16").Value = "Basic Aggregate Functions"

# Shade the new "Medium" difficulty cell orange (new fill, new cell style)
$ws.Range("B16").Interior.Color = 49407

# Turn the Link cell into a working hyperlink, then restore the plain
# "Hyperlink" cell style (Add() otherwise also bolds/underlines via a
# second, redundant style entry)
$ws.Hyperlinks.Add($ws.Range("E16"), "https://leetcode.com/problems/monthly-transactions-i/solutions/3871182/100-easy-fast-clean-solution/?envType=study-plan-v2&envId=top-sql-50 ")
$ws.Range("E16").Style = "Hyperlink"

# Extend the table (and therefore autofilter/dimension) to include the new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E16"))

# Move the active selection like the saved workbook shows
[void]$ws.Range("E26").Select()
